$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5162392258644104
$ws.Range("B1").Value = 3.123090505599976
$ws.Range("C1").Value = 6.275584697723389
$ws.Range("D1").Value = 2.803411245346069
$ws.Range("E1").Value = 1.907251119613647
